# Update the "status" column (C) for previously-written programs from the
# placeholder "<->" to their actual status (Done / rewrite / pending),
# matching the author's "Updated the filename for previously written
# programs" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C64:C72").Value = "Done"
$ws.Range("C73").Value = "rewrite"
$ws.Range("C74:C79").Value = "Done"
$ws.Range("C80").Value = "rewrite"
$ws.Range("C82").Value = "pending"
$ws.Range("C83:C84").Value = "Done"
$ws.Range("C85").Value = "pending"
$ws.Range("C86:C89").Value = "Done"
$ws.Range("C90:C91").Value = "pending"
$ws.Range("C92:C93").Value = "Done"
$ws.Range("C96:C100").Value = "Done"
$ws.Range("C101").Value = "rewrite"
$ws.Range("C102:C114").Value = "Done"

# Rows 113/114 also picked up the same row formatting already used by the
# surrounding "Done" rows (e.g. row 111/112) - copy that formatting down.
$ws.Range("A111:C111").Copy() | Out-Null
$ws.Range("A113:C114").PasteSpecial(-4122) | Out-Null

# Column B (the hyperlinked filenames) on rows 112-114 switched to the
# alternate hyperlink style already used elsewhere in the sheet (e.g. B87).
$ws.Range("B87").Copy() | Out-Null
$ws.Range("B112:B114").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Restore the cursor/selection and scroll position left by the editor.
$ws.Activate() | Out-Null
$ws.Range("B97").Select() | Out-Null
